{"js": "// Replace the date line and each two-digit-by-two-digit multiplication\n// prompt in the document body with its updated value. All \"old\" values are\n// unique exact strings within the document, so a scoped search + full-text\n// replace on each match is safe and preserves the existing run formatting\n// (font / size) since we replace in place rather than re-inserting text.\nconst replacements = [\n  [\"2025-06-16 Monday\", \"2025-06-17 Tuesday\"],\n  [\"21\u00d741=\", \"11\u00d714=\"],\n  [\"73\u00d726=\", \"87\u00d779=\"],\n  [\"70\u00d772=\", \"57\u00d724=\"],\n  [\"42\u00d781=\", \"59\u00d771=\"],\n  [\"91\u00d766=\", \"75\u00d777=\"],\n  [\"54\u00d783=\", \"72\u00d733=\"],\n  [\"79\u00d711=\", \"57\u00d748=\"],\n  [\"52\u00d778=\", \"78\u00d776=\"],\n  [\"83\u00d784=\", \"57\u00d753=\"],\n  [\"29\u00d768=\", \"15\u00d741=\"],\n  [\"26\u00d738=\", \"72\u00d763=\"],\n  [\"51\u00d767=\", \"37\u00d731=\"],\n  [\"62\u00d796=\", \"90\u00d798=\"],\n  [\"88\u00d799=\", \"33\u00d792=\"],\n  [\"36\u00d761=\", \"45\u00d746=\"],\n  [\"15\u00d745=\", \"47\u00d767=\"],\n  [\"79\u00d716=\", \"54\u00d745=\"],\n  [\"22\u00d794=\", \"51\u00d791=\"],\n  [\"78\u00d725=\", \"57\u00d711=\"],\n  [\"11\u00d720=\", \"26\u00d775=\"],\n  [\"33\u00d797=\", \"65\u00d765=\"],\n  [\"92\u00d764=\", \"49\u00d790=\"],\n  [\"13\u00d742=\", \"50\u00d714=\"],\n  [\"74\u00d711=\", \"38\u00d730=\"],\n  [\"32\u00d799=\", \"55\u00d776=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit-by-two-digit multiplication\n# prompt in the document body with its updated value. All \"old\" values\n# are unique exact strings within the document, so a global Find/Replace\n# per pair is safe and leaves formatting (font / size) untouched since\n# Word replaces the matched text in place.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-16 Monday\", \"2025-06-17 Tuesday\"),\n    @(\"21\u00d741=\", \"11\u00d714=\"),\n    @(\"73\u00d726=\", \"87\u00d779=\"),\n    @(\"70\u00d772=\", \"57\u00d724=\"),\n    @(\"42\u00d781=\", \"59\u00d771=\"),\n    @(\"91\u00d766=\", \"75\u00d777=\"),\n    @(\"54\u00d783=\", \"72\u00d733=\"),\n    @(\"79\u00d711=\", \"57\u00d748=\"),\n    @(\"52\u00d778=\", \"78\u00d776=\"),\n    @(\"83\u00d784=\", \"57\u00d753=\"),\n    @(\"29\u00d768=\", \"15\u00d741=\"),\n    @(\"26\u00d738=\", \"72\u00d763=\"),\n    @(\"51\u00d767=\", \"37\u00d731=\"),\n    @(\"62\u00d796=\", \"90\u00d798=\"),\n    @(\"88\u00d799=\", \"33\u00d792=\"),\n    @(\"36\u00d761=\", \"45\u00d746=\"),\n    @(\"15\u00d745=\", \"47\u00d767=\"),\n    @(\"79\u00d716=\", \"54\u00d745=\"),\n    @(\"22\u00d794=\", \"51\u00d791=\"),\n    @(\"78\u00d725=\", \"57\u00d711=\"),\n    @(\"11\u00d720=\", \"26\u00d775=\"),\n    @(\"33\u00d797=\", \"65\u00d765=\"),\n    @(\"92\u00d764=\", \"49\u00d790=\"),\n    @(\"13\u00d742=\", \"50\u00d714=\"),\n    @(\"74\u00d711=\", \"38\u00d730=\"),\n    @(\"32\u00d799=\", \"55\u00d776=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap = wdFindContinue\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace = wdReplaceAll\n    )\n}\n"}
